# Update cryptocurrency price/volume data to reflect the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. Numeric-looking "Price" values are
# prefixed with a literal apostrophe so Excel stores them as text (matching the
# original inline-string cells) instead of auto-converting them to numbers,
# which would silently drop things like trailing zeros (e.g. "2.60" -> 2.6).
$updates = [ordered]@{
    "D2" = "47.236.35"
    "E2" = "  +2.15%  "
    "D3" = "2.496.84"
    "E3" = "  +1.77%  "
    "E4" = "  +0.01%  "
    "D5" = "'323.47"
    "E5" = "  +0.73%  "
    "D6" = "'109.28"
    "E6" = "  +3.61%  "
    "E7" = "  +1.34%  "
    "D8" = "'0.999"
    "E8" = "  -0.04%  "
    "D9" = "'0.536"
    "E9" = "  -0.03%  "
    "D10" = "'39.21"
    "E10" = "  +8.59%  "
    "D11" = "'0.0814"
    "E11" = "  -0.15%  "
    "E12" = "  +0.70%  "
    "D13" = "'18.39"
    "E13" = "  +0.46%  "
    "D14" = "'7.21"
    "E14" = "  +1.66%  "
    "D15" = "2.887.06"
    "E15" = "  +1.47%  "
    "D16" = "2.494.74"
    "E16" = "  +2.67%  "
    "E17" = "  +1.23%  "
    "D18" = "47.172.43"
    "E18" = "  +2.27%  "
    "D19" = "'12.81"
    "E19" = "  +0.54%  "
    "E20" = "  +3.20%  "
    "E21" = "  +0.31%  "
    "D22" = "'2.71"
    "E22" = "  +13.07%  "
    "D23" = "'70.73"
    "E23" = "  -0.24%  "
    "D24" = "'247.32"
    "E24" = "  -0.30%  "
    "D25" = "'2.60"
    "E25" = "  +2.99%  "
    "D26" = "'26.07"
    "E26" = "  +0.27%  "
    "E27" = "  -0.02%  "
    "E28" = "  +2.20%  "
    "E29" = "  +3.43%  "
    "D30" = "'35.31"
    "E31" = "  +7.41%  "
    "D32" = "'49.84"
    "E32" = "  +1.08%  "
    "D33" = "'20.07"
    "E33" = "  +1.58%  "
    "D34" = "'5.43"
    "E34" = "  +1.50%  "
    "D35" = "'0.0789"
    "E35" = "  +2.88%  "
    "E36" = "  +0.19%  "
    "E37" = "  +4.24%  "
    "E38" = "  +2.44%  "
    "D39" = "'2.99"
    "E39" = "  +0.72%  "
    "E40" = "  +1.01%  "
    "D41" = "'121.76"
    "E41" = "  -2.46%  "
    "D42" = "'2.22"
    "E42" = "  -0.49%  "
    "D43" = "'21.37"
    "E43" = "  +2.68%  "
    "E44" = "  +2.13%  "
    "D45" = "1.994.63"
    "E45" = "  +0.88%  "
    "D46" = "'3.06"
    "E46" = "  +2.57%  "
    "E47" = "  -1.69%  "
    "D48" = "'1.79"
    "E48" = "  -3.53%  "
    "D49" = "'9.09"
    "E49" = "  -0.27%  "
    "E50" = "  +1.62%  "
    "D51" = "'56.74"
    "E51" = "  +4.16%  "
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
